$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Propagate the correct fill-colour ("style") into the rows that are new
#    or that need to change colour, by copying a single already-correctly
#    coloured cell across the destination row. Order matters: grab the
#    colour that a row currently has *before* that row's own content/colour
#    is overwritten by a later step.
# ---------------------------------------------------------------------------

# New row 9 ("categories/image" boundary case) reuses the colour that (old)
# row 8 currently has (the blue fill) - copy it before row 8 is recoloured.
$ws.Range("A8").Copy($ws.Range("A9:G9"))

# Row 8 becomes the "comments len = 0" boundary case and must switch to the
# colour that (old) row 7 currently has (the purple fill) - copy it before
# row 7 is recoloured.
$ws.Range("A7").Copy($ws.Range("A8:G8"))

# Row 7 becomes the "servings < 0" boundary case and must switch to the same
# colour as row 6 (orange fill).
$ws.Range("A6").Copy($ws.Range("A7:G7"))

# New row 5 ("preparationTime < 0" boundary case) uses the same colour as
# row 4 (yellow fill).
$ws.Range("A4").Copy($ws.Range("A5:G5"))

# ---------------------------------------------------------------------------
# 2) Now write the actual cell contents for every row in the table.
# ---------------------------------------------------------------------------

# Row 4 - preparationTime "equals(0)" boundary
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "len > 0"
$ws.Range("C4").Value = "equals(0)"
$ws.Range("D4").Value = "> 0"
$ws.Range("E4").Value = "len > 0"
$ws.Range("F4").Value = "X"
$ws.Range("G4").Value = "null"

# Row 5 (new) - preparationTime "< 0" boundary
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "len > 0"
$ws.Range("C5").Value = "< 0"
$ws.Range("D5").Value = "> 0"
$ws.Range("E5").Value = "len > 0"
$ws.Range("F5").Value = "X"
$ws.Range("G5").Value = "null"

# Row 6 - servings "equals(0)" boundary
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "len > 0"
$ws.Range("C6").Value = "> 0"
$ws.Range("D6").Value = "equals(0)"
$ws.Range("E6").Value = "len > 0"
$ws.Range("F6").Value = "X"
$ws.Range("G6").Value = "null"

# Row 7 - servings "< 0" boundary
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "len > 0"
$ws.Range("C7").Value = "> 0"
$ws.Range("D7").Value = "< 0"
$ws.Range("E7").Value = "len > 0"
$ws.Range("F7").Value = "X"
$ws.Range("G7").Value = "null"

# Row 8 - comments "len = 0" boundary
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "len > 0"
$ws.Range("C8").Value = "> 0"
$ws.Range("D8").Value = "> 0"
$ws.Range("E8").Value = "len = 0"
$ws.Range("F8").Value = "X"
$ws.Range("G8").Value = "null"

# Row 9 (new) - categories/image boundary
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "len > 0"
$ws.Range("C9").Value = "> 0"
$ws.Range("D9").Value = "> 0"
$ws.Range("E9").Value = "len > 0"
$ws.Range("F9").Value = "X"
$ws.Range("G9").Value = "image"

# ---------------------------------------------------------------------------
# 3) Update the "Number of States" summary table lower on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = 3
$ws.Range("B15").Value = 3

# ---------------------------------------------------------------------------
# 4) Match the final selection left behind on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("H30").Select()
